# Append the new training-log row (row 26) to the "training" sheet,
# matching the existing log-row layout (columns A:N).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("training")

$row = 26

$ws.Cells.Item($row, 1).Value = "2025-09-16 14:27:00"
$ws.Cells.Item($row, 2).Value = "training"
$ws.Cells.Item($row, 3).Value = "configs/training/2025-09-16/x/0001"

# "note" column is blank for this run, but (like every other row) the cell
# itself is still present in the sheet as an empty text cell rather than
# being entirely absent. A bare value of "" doesn't materialize the cell,
# so write a lone text-prefix apostrophe (-> empty text) and then strip the
# quote-prefix formatting it introduces so the cell matches its neighbours.
$ws.Cells.Item($row, 4).Value = "'"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = "['cross_entropy', 'spectral_entropy']"
$ws.Cells.Item($row, 6).Value = "[1.0, 0.01]"
$ws.Cells.Item($row, 7).Value = "['torch.optim.adamw.AdamW', 'torch.optim.adam.Adam']"
$ws.Cells.Item($row, 8).Value = "[0.001, 0.001]"
$ws.Cells.Item($row, 9).Value = 128
$ws.Cells.Item($row, 10).Value = 128
$ws.Cells.Item($row, 11).Value = "general_utils.ml.training.NoImprovementStopping"
$ws.Cells.Item($row, 12).Value = 20
$ws.Cells.Item($row, 13).Value = 0.00001
$ws.Cells.Item($row, 14).Value = 500
